$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.384.64'
$ws.Range("E2").Value = '  +9.60%  '
$ws.Range("D3").Value = '3.238.06'
$ws.Range("E3").Value = '  +4.79%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '400.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +8.94%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.557'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.44%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.627'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.81%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.80'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0901'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.48%  '
$ws.Range("E12").Value = '  +2.20%  '
$ws.Range("D13").Value = '3.754.23'
$ws.Range("E13").Value = '  +4.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.25'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.09'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.65%  '
$ws.Range("E16").Value = '  +7.79%  '
$ws.Range("D17").Value = '3.242.44'
$ws.Range("E17").Value = '  +4.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.68'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.41%  '
$ws.Range("D19").Value = '56.252.10'
$ws.Range("E19").Value = '  +9.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.90%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000103'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.25%  '
$ws.Range("B22").Value = 'InternetComputer(DFINITY)'
$ws.Range("C22").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.21'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '290.58'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.84%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.24'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.16'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.17%  '
$ws.Range("E29").Value = '  +3.81%  '
$ws.Range("E30").Value = '  -0.06%  '
$ws.Range("E31").Value = '  +5.33%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +9.92%  '
$ws.Range("B33").Value = 'VeChain'
$ws.Range("C33").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0499'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.49%  '
$ws.Range("B34").Value = 'InjectiveProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '37.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.40%  '
$ws.Range("E35").Value = '  +2.51%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.52'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.59'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.37%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.13'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +25.01%  '
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '138.68'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.79%  '
$ws.Range("E41").Value = '  +3.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.03'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.71%  '
$ws.Range("B43").Value = 'Celestia'
$ws.Range("C43").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.13'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.73%  '
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.286'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.59%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.119'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.09%  '
$ws.Range("E47").Value = '  +41.43%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '2.146.82'
$ws.Range("E48").Value = '  +4.66%  '
$ws.Range("B49").Value = 'WEMIXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.09'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.97%  '
$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.44'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0351'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +7.61%  '